$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17; existing rows 17-25 shift down to 18-26.
$ws.Rows.Item(17).Insert()

# Populate the new row 17 with this week's price entry (weekly update).
$ws.Range("A17").Value = 1
$ws.Range("B17").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C17").Value = "Arica y Parinacota"
$ws.Range("D17").Value = 44762
$ws.Range("E17").Value = 15
$ws.Range("F17").Value = 100112043
$ws.Range("G17").Value = "Pepino dulce"
$ws.Range("H17").Value = "Cultivar IV Región"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 160
$ws.Range("K17").Value = 15000
$ws.Range("L17").Value = 16000
$ws.Range("M17").Value = 15500
$ws.Range("N17").Value = "$/bandeja 18 kilos"
$ws.Range("O17").Value = "Provincia de Limarí"
$ws.Range("P17").Value = 861
$ws.Range("Q17").Value = 18
$ws.Range("R17").Value = "Hortaliza"
